# Auto-generated script to apply scheduled market-data refresh to Titan_Profits leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1018.05554
$ws.Range("I43").Value = 941.2857
$ws.Range("J43").Value = 1066.909
$ws.Range("K43").Value = 941.2857
$ws.Range("L43").Value = 1066.909
$ws.Range("M43").Value = -872.2857
$ws.Range("N43").Value = -1204.909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3986.1667
$ws.Range("I58").Value = 107
$ws.Range("J58").Value = 6757
$ws.Range("K58").Value = 321
$ws.Range("L58").Value = 20271
$ws.Range("M58").Value = -171
$ws.Range("N58").Value = -20571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2381
$ws.Range("I70").Value = 1624.75
$ws.Range("J70").Value = 2717.111
$ws.Range("K70").Value = 4874.25
$ws.Range("L70").Value = 8151.333
$ws.Range("M70").Value = -4604.25
$ws.Range("N70").Value = -8691.332999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2381
$ws.Range("I73").Value = 1624.75
$ws.Range("J73").Value = 2717.111
$ws.Range("K73").Value = 4874.25
$ws.Range("L73").Value = 8151.333
$ws.Range("M73").Value = -3938.25
$ws.Range("N73").Value = -10023.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 6536789.5
$ws.Range("I92").Value = 7937378
$ws.Range("J92").Value = 709.3333
$ws.Range("K92").Value = 7937378
$ws.Range("L92").Value = 709.3333
$ws.Range("M92").Value = -7936130
$ws.Range("N92").Value = -3205.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 11364310
$ws.Range("I99").Value = 18182100
$ws.Range("J99").Value = 1326.6666
$ws.Range("K99").Value = 54546300
$ws.Range("L99").Value = 3979.9998
$ws.Range("M99").Value = -54544802
$ws.Range("N99").Value = -6975.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7577230.5
$ws.Range("I100").Value = 11112451
$ws.Range("K100").Value = 11112451
$ws.Range("M100").Value = -11111910

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 34054.47
$ws.Range("I132").Value = 35133.645
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 105400.935
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -102870.935
$ws.Range("N132").Value = -6860

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 25002014
$ws.Range("I137").Value = 45455700
$ws.Range("J137").Value = 3061.611
$ws.Range("K137").Value = 136367100
$ws.Range("L137").Value = 9184.832999999999
$ws.Range("M137").Value = -136364550
$ws.Range("N137").Value = -14284.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21881.314
$ws.Range("I32").Value = 5705.484
$ws.Range("J32").Value = 113054.18
$ws.Range("K32").Value = 5705.484
$ws.Range("L32").Value = 113054.18
$ws.Range("M32").Value = -5418.484
$ws.Range("N32").Value = -113628.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7046.591
$ws.Range("I74").Value = 1510.7646
$ws.Range("J74").Value = 25868.4
$ws.Range("K74").Value = 1510.7646
$ws.Range("L74").Value = 25868.4
$ws.Range("M74").Value = -636.7646
$ws.Range("N74").Value = -27616.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7046.591
$ws.Range("I77").Value = 1510.7646
$ws.Range("J77").Value = 25868.4
$ws.Range("K77").Value = 7553.823
$ws.Range("L77").Value = 129342
$ws.Range("M77").Value = -3185.823
$ws.Range("N77").Value = -138078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2719.15
$ws.Range("I102").Value = 2755.1428
$ws.Range("J102").Value = 2635.1667
$ws.Range("K102").Value = 2755.1428
$ws.Range("L102").Value = 2635.1667
$ws.Range("M102").Value = -1133.1428
$ws.Range("N102").Value = -5879.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2471.889
$ws.Range("I132").Value = 1954.4318
$ws.Range("J132").Value = 4748.7
$ws.Range("K132").Value = 5863.2954
$ws.Range("L132").Value = 14246.1
$ws.Range("M132").Value = -3333.2954
$ws.Range("N132").Value = -19306.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 59000
$ws.Range("J139").Value = 59000
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1366
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 985.1053000000001
$ws.Range("J64").Value = 474.46667
$ws.Range("L64").Value = 474.46667
$ws.Range("N64").Value = -924.46667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 985.1053000000001
$ws.Range("J67").Value = 474.46667
$ws.Range("L67").Value = 474.46667
$ws.Range("N67").Value = -2034.46667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2079.2856
$ws.Range("I107").Value = 2079.2856
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2079.2856
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -159.2856000000002
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 51000
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3494.3333
$ws.Range("I134").Value = 2173.55
$ws.Range("K134").Value = 6520.650000000001
$ws.Range("M134").Value = -3985.650000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3833.9395
$ws.Range("I31").Value = 1745.1666
$ws.Range("J31").Value = 5574.5835
$ws.Range("K31").Value = 1745.1666
$ws.Range("L31").Value = 5574.5835
$ws.Range("M31").Value = -1450.1666
$ws.Range("N31").Value = -6164.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3833.9395
$ws.Range("I34").Value = 1745.1666
$ws.Range("J34").Value = 5574.5835
$ws.Range("K34").Value = 1745.1666
$ws.Range("L34").Value = 5574.5835
$ws.Range("M34").Value = -1543.1666
$ws.Range("N34").Value = -5978.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 33334804
$ws.Range("I58").Value = 37038136
$ws.Range("K58").Value = 37038136
$ws.Range("M58").Value = -37037933

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 370.08823
$ws.Range("I107").Value = 303.7037
$ws.Range("J107").Value = 626.1429000000001
$ws.Range("K107").Value = 303.7037
$ws.Range("L107").Value = 626.1429000000001
$ws.Range("M107").Value = 1616.2963
$ws.Range("N107").Value = -4466.1429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 28848460
$ws.Range("I134").Value = 30304412
$ws.Range("J134").Value = 26319704
$ws.Range("K134").Value = 90913236
$ws.Range("L134").Value = 78959112
$ws.Range("M134").Value = -90910701
$ws.Range("N134").Value = -78964182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 33334804
$ws.Range("I136").Value = 37038136
$ws.Range("K136").Value = 111114408
$ws.Range("M136").Value = -111111858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5348351.5
$ws.Range("I113").Value = 630
$ws.Range("J113").Value = 6494292
$ws.Range("K113").Value = 1890
$ws.Range("L113").Value = 19482876
$ws.Range("M113").Value = 280
$ws.Range("N113").Value = -19487216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5849439
$ws.Range("I131").Value = 588.8889
$ws.Range("J131").Value = 6946098.5
$ws.Range("K131").Value = 1766.6667
$ws.Range("L131").Value = 20838295.5
$ws.Range("M131").Value = 3273.3333
$ws.Range("N131").Value = -20848375.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8750
$ws.Range("J92").Value = 8750
$ws.Range("L92").Value = 8750
$ws.Range("N92").Value = -12494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1109.35
$ws.Range("I97").Value = 1241.1428
$ws.Range("J97").Value = 801.8333
$ws.Range("K97").Value = 1241.1428
$ws.Range("L97").Value = 801.8333
$ws.Range("M97").Value = -745.1428000000001
$ws.Range("N97").Value = -1793.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3017.6775
$ws.Range("I132").Value = 2806.524
$ws.Range("J132").Value = 3461.1
$ws.Range("K132").Value = 8419.572
$ws.Range("L132").Value = 10383.3
$ws.Range("M132").Value = -5889.572
$ws.Range("N132").Value = -15443.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8333689.5
$ws.Range("I16").Value = 9091270
$ws.Range("J16").Value = 302
$ws.Range("K16").Value = 9091270
$ws.Range("L16").Value = 302
$ws.Range("M16").Value = -9091100
$ws.Range("N16").Value = -642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1840767.4
$ws.Range("I100").Value = 6252070.5
$ws.Range("J100").Value = 2724.2917
$ws.Range("K100").Value = 6252070.5
$ws.Range("L100").Value = 2724.2917
$ws.Range("M100").Value = -6251529.5
$ws.Range("N100").Value = -3806.2917

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 744829.75
$ws.Range("I81").Value = 2224511.5
$ws.Range("K81").Value = 4449023
$ws.Range("M81").Value = -4447962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 744829.75
$ws.Range("I84").Value = 2224511.5
$ws.Range("K84").Value = 22245115
$ws.Range("M84").Value = -22239811

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 59432.41
$ws.Range("I126").Value = 62869.188
$ws.Range("J126").Value = 4444
$ws.Range("K126").Value = 188607.564
$ws.Range("L126").Value = 13332
$ws.Range("M126").Value = -186137.564
$ws.Range("N126").Value = -18272
